$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,14

$data[0,0] = 1.261954552005932
$data[0,1] = 0.174196266374679
$data[0,2] = 0.07859522710324995
$data[0,3] = 0.07376806398248448
$data[0,4] = 0
$data[0,5] = 0.00243886799639053
$data[0,6] = 0
$data[0,7] = 0.7385653484599146
$data[0,8] = 0
$data[0,9] = 0
$data[0,10] = 0.235314674209647
$data[0,11] = 0
$data[0,12] = 1.180148526259536
$data[0,13] = 3.226305097233421
$data[1,0] = 1.154455080356456
$data[1,1] = 0.1566964042896473
$data[1,2] = 0.07127836230610285
$data[1,3] = 0.07423817962832402
$data[1,4] = 0
$data[1,5] = 0.002442229259409349
$data[1,6] = 0
$data[1,7] = 0.742670473319933
$data[1,8] = 0
$data[1,9] = 0
$data[1,10] = 0.2255750777908361
$data[1,11] = 0
$data[1,12] = 1.19292675491992
$data[1,13] = 3.210457245889131
$data[2,0] = 1.088719809577412
$data[2,1] = 0.1458786228336351
$data[2,2] = 0.06682204113344881
$data[2,3] = 0.07456088428630991
$data[2,4] = 0
$data[2,5] = 0.002444404158545892
$data[2,6] = 0
$data[2,7] = 0.7456702100432402
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0.2197199117820077
$data[2,11] = 0
$data[2,12] = 1.201268713229638
$data[2,13] = 3.202968835800164
$data[3,0] = 1.062001216078272
$data[3,1] = 0.1414520916558217
$data[3,2] = 0.06501512982498525
$data[3,3] = 0.07470095547736477
$data[3,4] = 0
$data[3,5] = 0.002445318463659457
$data[3,6] = 0
$data[3,7] = 0.7470129307300404
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0.2173653471737396
$data[3,11] = 0
$data[3,12] = 1.204792880396539
$data[3,13] = 3.200480047049552
$data[4,0] = 1.057568821012069
$data[4,1] = 0.1407159736427275
$data[4,2] = 0.06471564048375456
$data[4,3] = 0.07472473169665861
$data[4,4] = 0
$data[4,5] = 0.002445471977745662
$data[4,6] = 0
$data[4,7] = 0.7472431490194396
$data[4,8] = 0
$data[4,9] = 0
$data[4,10] = 0.2169762748467576
$data[4,11] = 0
$data[4,12] = 1.20538560072335
$data[4,13] = 3.200100747266873
$data[5,0] = 1.088359192042162
$data[5,1] = 0.1458189986223886
$data[5,2] = 0.06679763580442
$data[5,3] = 0.07456273864498897
$data[5,4] = 0
$data[5,5] = 0.002444416375825673
$data[5,6] = 0
$data[5,7] = 0.7456878315913968
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0.219688029860265
$data[5,11] = 0
$data[5,12] = 1.201315736254433
$data[5,13] = 3.202932993745321
$data[6,0] = 1.224833418626702
$data[6,1] = 0.1681774820262092
$data[6,2] = 0.07606481752209504
$data[6,3] = 0.07392309466344216
$data[6,4] = 0
$data[6,5] = 0.002440003959349593
$data[6,6] = 0
$data[6,7] = 0.7398812285058014
$data[6,8] = 0
$data[6,9] = 0
$data[6,10] = 0.2319305375421692
$data[6,11] = 0
$data[6,12] = 1.184451504636463
$data[6,13] = 3.22037474737229
$data[7,0] = 1.494560939655912
$data[7,1] = 0.2114420842765696
$data[7,2] = 0.09452836721554547
$data[7,3] = 0.07293879798887559
$data[7,4] = 0
$data[7,5] = 0.002432228578048645
$data[7,6] = 0
$data[7,7] = 0.732305584107344
$data[7,8] = 0
$data[7,9] = 0
$data[7,10] = 0.2569297870633562
$data[7,11] = 0
$data[7,12] = 1.155316413127103
$data[7,13] = 3.272422403133305
$data[8,0] = 1.693979431654043
$data[8,1] = 0.2428741305002688
$data[8,2] = 0.1082761172732631
$data[8,3] = 0.07238010914301185
$data[8,4] = 0
$data[8,5] = 0.002427045350800751
$data[8,6] = 0
$data[8,7] = 0.7290764498606137
$data[8,8] = 0
$data[8,9] = 0
$data[8,10] = 0.2759036320956483
$data[8,11] = 0
$data[8,12] = 1.136309294854577
$data[8,13] = 3.321623982593735
$data[9,0] = 1.78496605701946
$data[9,1] = 0.2570964597386194
$data[9,2] = 0.114571209610375
$data[9,3] = 0.07216163318665103
$data[9,4] = 0
$data[9,5] = 0.002424801135704844
$data[9,6] = 0
$data[9,7] = 0.7281177684424094
$data[9,8] = 0
$data[9,9] = 0
$data[9,10] = 0.2846677490309588
$data[9,11] = 0
$data[9,12] = 1.128183187027616
$data[9,13] = 3.34640608923911
$data[10,0] = 1.819458259749467
$data[10,1] = 0.2624710531641483
$data[10,2] = 0.1169609857892482
$data[10,3] = 0.0720840294532632
$data[10,4] = 0
$data[10,5] = 0.002423967564766307
$data[10,6] = 0
$data[10,7] = 0.7278283459446513
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0.2880055945948214
$data[10,11] = 0
$data[10,12] = 1.125180886372206
$data[10,13] = 3.356136882399426
$data[11,0] = 1.812028097755444
$data[11,1] = 0.261314033729434
$data[11,2] = 0.116446038950059
$data[11,3] = 0.07210051474469914
$data[11,4] = 0
$data[11,5] = 0.00242414636704906
$data[11,6] = 0
$data[11,7] = 0.7278874007367335
$data[11,8] = 0
$data[11,9] = 0
$data[11,10] = 0.2872858824720907
$data[11,11] = 0
$data[11,12] = 1.125824154249436
$data[11,13] = 3.35402575861491
$data[12,0] = 1.787803006346849
$data[12,1] = 0.2575388538403161
$data[12,2] = 0.1147676982267143
$data[12,3] = 0.07215514590942185
$data[12,4] = 0
$data[12,5] = 0.002424732231634197
$data[12,6] = 0
$data[12,7] = 0.7280924808987805
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0.2849419737068359
$data[12,11] = 0
$data[12,12] = 1.127934684888999
$data[12,13] = 3.347199699127998
$data[13,0] = 1.772969288865227
$data[13,1] = 0.2552249986276536
$data[13,2] = 0.1137404436230725
$data[13,3] = 0.07218927689461196
$data[13,4] = 0
$data[13,5] = 0.002425093206388058
$data[13,6] = 0
$data[13,7] = 0.7282276914699963
$data[13,8] = 0
$data[13,9] = 0
$data[13,10] = 0.2835087436721864
$data[13,11] = 0
$data[13,12] = 1.129237199240926
$data[13,13] = 3.343063681526132
$data[14,0] = 1.688038566149771
$data[14,1] = 0.2419431222930086
$data[14,2] = 0.107865550166423
$data[14,3] = 0.07239510478607514
$data[14,4] = 0
$data[14,5] = 0.002427194295920709
$data[14,6] = 0
$data[14,7] = 0.7291493897799484
$data[14,8] = 0
$data[14,9] = 0
$data[14,10] = 0.2753335461476638
$data[14,11] = 0
$data[14,12] = 1.1368508302156
$data[14,13] = 3.320052813340681
$data[15,0] = 1.636004580616088
$data[15,1] = 0.2337754971831032
$data[15,2] = 0.1042720675800979
$data[15,3] = 0.07253050947403317
$data[15,4] = 0
$data[15,5] = 0.002428512301110337
$data[15,6] = 0
$data[15,7] = 0.7298456795739767
$data[15,8] = 0
$data[15,9] = 0
$data[15,10] = 0.2703523208324583
$data[15,11] = 0
$data[15,12] = 1.141654853552893
$data[15,13] = 3.306551972324371
$data[16,0] = 1.606101514100374
$data[16,1] = 0.2290705242634772
$data[16,2] = 0.1022090643090792
$data[16,3] = 0.07261174860445685
$data[16,4] = 0
$data[16,5] = 0.002429281085597146
$data[16,6] = 0
$data[16,7] = 0.7302941818712299
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0.2674997554498617
$data[16,11] = 0
$data[16,12] = 1.144466972415728
$data[16,13] = 3.299012477543357
$data[17,0] = 1.595981270094171
$data[17,1] = 0.2274762735312663
$data[17,2] = 0.1015112307683665
$data[17,3] = 0.07263983155498899
$data[17,4] = 0
$data[17,5] = 0.002429543223554661
$data[17,6] = 0
$data[17,7] = 0.7304542765379054
$data[17,8] = 0
$data[17,9] = 0
$data[17,10] = 0.2665360744014293
$data[17,11] = 0
$data[17,12] = 1.145427516810173
$data[17,13] = 3.29649848385634
$data[18,0] = 1.64154105693359
$data[18,1] = 0.2346456979128959
$data[18,2] = 0.1046541988129945
$data[18,3] = 0.07251574789184012
$data[18,4] = 0
$data[18,5] = 0.002428370890070884
$data[18,6] = 0
$data[18,7] = 0.7297665872911097
$data[18,8] = 0
$data[18,9] = 0
$data[18,10] = 0.2708812869409627
$data[18,11] = 0
$data[18,12] = 1.141138388274229
$data[18,13] = 3.307965778297586
$data[19,0] = 1.79491749440092
$data[19,1] = 0.2586480177828321
$data[19,2] = 0.115260505468072
$data[19,3] = 0.07213896024185829
$data[19,4] = 0
$data[19,5] = 0.002424559708342655
$data[19,6] = 0
$data[19,7] = 0.7280302442831541
$data[19,8] = 0
$data[19,9] = 0
$data[19,10] = 0.2856299193926759
$data[19,11] = 0
$data[19,12] = 1.127312738170552
$data[19,13] = 3.349195268751316
$data[20,0] = 1.895376087069508
$data[20,1] = 0.2742701985748965
$data[20,2] = 0.1222271059957762
$data[20,3] = 0.07192259909741949
$data[20,4] = 0
$data[20,5] = 0.00242216364916248
$data[20,6] = 0
$data[20,7] = 0.7273246042322512
$data[20,8] = 0
$data[20,9] = 0
$data[20,10] = 0.2953801110308945
$data[20,11] = 0
$data[20,12] = 1.118713365220408
$data[20,13] = 3.378160469920601
$data[21,0] = 1.841739907642761
$data[21,1] = 0.2659383106092719
$data[21,2] = 0.1185057056499517
$data[21,3] = 0.07203534055683036
$data[21,4] = 0
$data[21,5] = 0.002423433825729162
$data[21,6] = 0
$data[21,7] = 0.7276618720770429
$data[21,8] = 0
$data[21,9] = 0
$data[21,10] = 0.2901660952444161
$data[21,11] = 0
$data[21,12] = 1.123263053196077
$data[21,13] = 3.362516022288048
$data[22,0] = 1.639037978621673
$data[22,1] = 0.234252309114737
$data[22,2] = 0.1044814281370634
$data[22,3] = 0.0725224110338587
$data[22,4] = 0
$data[22,5] = 0.002428434787660271
$data[22,6] = 0
$data[22,7] = 0.7298021947978839
$data[22,8] = 0
$data[22,9] = 0
$data[22,10] = 0.2706421064535505
$data[22,11] = 0
$data[22,12] = 1.141371725786875
$data[22,13] = 3.307325904185177
$data[23,0] = 1.421370773247077
$data[23,1] = 0.1998000626224723
$data[23,2] = 0.08950181423456627
$data[23,3] = 0.07317618218307054
$data[23,4] = 0
$data[23,5] = 0.002434238670099946
$data[23,6] = 0
$data[23,7] = 0.73394561666219
$data[23,8] = 0
$data[23,9] = 0
$data[23,10] = 0.2500604186415103
$data[23,11] = 0
$data[23,12] = 1.162776998640787
$data[23,13] = 3.25642312388814

$ws.Range("B2:O25").Value = $data
